$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added at the top of the Pomelo data
# block (row 199). Insert a blank row there, which shifts the existing
# rows 199-248 down to 200-249 (carrying their values/formatting with
# them), then populate the new row 199 with the latest observation.
$ws.Rows("199:199").Insert()

# Columns that stay constant for every record in this block (market,
# region, product taxonomy, unit, origin, kg/unit) are copied from the
# row immediately below (the record that used to be row 199).
$ws.Range("A199").Value = $ws.Range("A200").Value()
$ws.Range("B199").Value = $ws.Range("B200").Value()
$ws.Range("C199").Value = $ws.Range("C200").Value()
$ws.Range("D199").Value = 44641
$ws.Range("E199").Value = $ws.Range("E200").Value()
$ws.Range("F199").Value = $ws.Range("F200").Value()
$ws.Range("G199").Value = $ws.Range("G200").Value()
$ws.Range("H199").Value = $ws.Range("H200").Value()
$ws.Range("I199").Value = $ws.Range("I200").Value()
$ws.Range("J199").Value = $ws.Range("J200").Value()
$ws.Range("K199").Value = "Start Ruby"
$ws.Range("L199").Value = "Primera"
$ws.Range("M199").Value = 80
$ws.Range("N199").Value = 12000
$ws.Range("O199").Value = 13000
$ws.Range("P199").Value = 12500
$ws.Range("Q199").Value = $ws.Range("Q200").Value()
$ws.Range("R199").Value = $ws.Range("R200").Value()
$ws.Range("S199").Value = 893
$ws.Range("T199").Value = $ws.Range("T200").Value()
